$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.157.44"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3
$ws.Range("D3").Value = "2.962.86"
$ws.Range("E3").Value = "  +1.24%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'380.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.92%  "

# Row 6
$ws.Range("D6").Value = "'102.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.48%  "

# Row 7
$ws.Range("E7").Value = "  +2.67%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  +2.17%  "

# Row 10
$ws.Range("D10").Value = "'36.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.00%  "

# Row 11
$ws.Range("E11").Value = "  -0.72%  "

# Row 12
$ws.Range("E12").Value = "  +2.21%  "

# Row 13
$ws.Range("D13").Value = "3.426.53"
$ws.Range("E13").Value = "  +1.17%  "

# Row 14
$ws.Range("D14").Value = "'7.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.13%  "

# Row 15
$ws.Range("D15").Value = "'18.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.42%  "

# Row 16
$ws.Range("D16").Value = "2.969.47"
$ws.Range("E16").Value = "  +1.88%  "

# Row 17
$ws.Range("D17").Value = "'11.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18
$ws.Range("E18").Value = "  +3.58%  "

# Row 19
$ws.Range("D19").Value = "51.205.96"
$ws.Range("E19").Value = "  +0.55%  "

# Row 20
$ws.Range("D20").Value = "'3.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.36%  "

# Row 21
$ws.Range("D21").Value = "'12.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.95%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  +1.15%  "

# Row 23
$ws.Range("D23").Value = "'70.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.78%  "

# Row 24
$ws.Range("D24").Value = "'266.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.41%  "

# Row 25
$ws.Range("D25").Value = "'3.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.14%  "

# Row 26
$ws.Range("D26").Value = "'7.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.58%  "

# Row 27
$ws.Range("D27").Value = "'7.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.39%  "

# Row 28
$ws.Range("E28").Value = "  -0.01%  "

# Row 29
$ws.Range("D29").Value = "'25.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.32%  "

# Row 30
$ws.Range("E30").Value = "  +2.20%  "

# Row 31
$ws.Range("E31").Value = "  -0.56%  "

# Row 32
$ws.Range("D32").Value = "'10.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.04%  "

# Row 33
$ws.Range("D33").Value = "'34.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.87%  "

# Row 34
$ws.Range("D34").Value = "'51.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.46%  "

# Row 35
$ws.Range("D35").Value = "'2.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.60%  "

# Row 36
$ws.Range("D36").Value = "'0.0437"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("D38").Value = "'3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.76%  "

# Row 39
$ws.Range("E39").Value = "  +2.00%  "

# Row 40
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'16.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.19%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.94%  "

# Row 42
$ws.Range("E42").Value = "  +4.01%  "

# Row 43
$ws.Range("D43").Value = "'124.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.73%  "

# Row 44
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.97%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'21.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.22%  "

# Row 46
$ws.Range("E46").Value = "  +0.09%  "

# Row 47
$ws.Range("D47").Value = "'2.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.19%  "

# Row 48
$ws.Range("E48").Value = "  -0.56%  "

# Row 49
$ws.Range("D49").Value = "2.034.21"
$ws.Range("E49").Value = "  +3.13%  "

# Row 50
$ws.Range("D50").Value = "'0.0326"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "

# Row 51
$ws.Range("D51").Value = "'0.512"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.84%  "
